$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2: new labels above each CPT (Conditional Probability Table) block ----
$ws.Range("B2").Value = "Placebo"
$ws.Range("F2").Value = "Medication "
$ws.Range("J2").Value = "No treatment"

# ---- Row 4/5: update probability values in the first three tables ----
$ws.Range("C4").Value = 0.09
$ws.Range("D4").Value = 0.91
$ws.Range("G4").Value = 0.072
$ws.Range("H4").Value = 0.928
$ws.Range("K4").Value = 0.043
$ws.Range("L4").Value = 0.957

$ws.Range("C5").Value = 0.902
$ws.Range("D5").Value = 0.098
$ws.Range("G5").Value = 0.902
$ws.Range("H5").Value = 0.098
$ws.Range("K5").Value = 0.39
$ws.Range("L5").Value = 0.61

# ---- Row 6: new labels ----
$ws.Range("B6").Value = "Survive"
$ws.Range("F6").Value = "Healthy"

# ---- Row 7: F7 header text change, clear J7:L7 content & remove border ----
$ws.Range("F7").Value = "Z"
$ws.Range("J7:L7").ClearContents()
$ws.Range("J7:L7").Borders.LineStyle = -4142  # xlLineStyleNone

# ---- Row 8: remove border formatting from J8:L8 (keep them blank) ----
$ws.Range("J8:L8").Borders.LineStyle = -4142  # xlLineStyleNone

# ---- Row 11/12: add TT/TF/FT/FF mini legend ----
$ws.Range("G11").Value = "TT"
$ws.Range("H11").Value = "TF"
$ws.Range("G12").Value = "FT"
$ws.Range("H12").Value = "FF"

# ---- Row 18/19: convert plain text URLs into real hyperlinks ----
$ws.Hyperlinks.Add($ws.Range("B19"), "https://medium.com/analytics-vidhya/comparison-of-two-data-sets-using-python-a24a6d8beb13")
$ws.Hyperlinks.Add($ws.Range("B18"), "https://stackoverflow.com/questions/55926173/compare-elements-in-dataframe-columns-for-each-row-python")

# ---- Selection moved ----
$ws.Range("H16").Select() | Out-Null
